$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update Testable column (B) values from "y" to "n" for rows 2 and 3
$ws.Range("B2").Value = "n"
$ws.Range("B3").Value = "n"

# Update the view: scroll back to A1 (remove topLeftCell freeze at D1) and select B7
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B7").Select()
